{"js": "// Update the worksheet date and the two-digit multiplication problems\n// to the new set of values, matching the target revision.\nconst replacements = [\n  [\"2025-10-18 Saturday\", \"2025-10-19 Sunday\"],\n  [\"59\u00d780=\", \"87\u00d782=\"],\n  [\"93\u00d743=\", \"71\u00d773=\"],\n  [\"96\u00d794=\", \"90\u00d729=\"],\n  [\"93\u00d735=\", \"12\u00d755=\"],\n  [\"34\u00d750=\", \"70\u00d788=\"],\n  [\"87\u00d756=\", \"42\u00d784=\"],\n  [\"73\u00d794=\", \"73\u00d730=\"],\n  [\"37\u00d713=\", \"94\u00d727=\"],\n  [\"60\u00d712=\", \"27\u00d791=\"],\n  [\"83\u00d723=\", \"26\u00d752=\"],\n  [\"16\u00d723=\", \"15\u00d723=\"],\n  [\"11\u00d796=\", \"70\u00d798=\"],\n  [\"38\u00d750=\", \"45\u00d739=\"],\n  [\"83\u00d721=\", \"58\u00d798=\"],\n  [\"47\u00d772=\", \"48\u00d785=\"],\n  [\"74\u00d721=\", \"65\u00d732=\"],\n  [\"99\u00d776=\", \"84\u00d739=\"],\n  [\"58\u00d742=\", \"93\u00d720=\"],\n  [\"93\u00d775=\", \"86\u00d744=\"],\n  [\"30\u00d740=\", \"31\u00d755=\"],\n  [\"78\u00d793=\", \"12\u00d778=\"],\n  [\"90\u00d756=\", \"71\u00d724=\"],\n  [\"81\u00d724=\", \"53\u00d758=\"],\n  [\"58\u00d713=\", \"50\u00d761=\"],\n  [\"88\u00d770=\", \"85\u00d760=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the two-digit multiplication problems\n# to the new set of values, matching the target revision.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-18 Saturday\", \"2025-10-19 Sunday\"),\n    @(\"59\u00d780=\", \"87\u00d782=\"),\n    @(\"93\u00d743=\", \"71\u00d773=\"),\n    @(\"96\u00d794=\", \"90\u00d729=\"),\n    @(\"93\u00d735=\", \"12\u00d755=\"),\n    @(\"34\u00d750=\", \"70\u00d788=\"),\n    @(\"87\u00d756=\", \"42\u00d784=\"),\n    @(\"73\u00d794=\", \"73\u00d730=\"),\n    @(\"37\u00d713=\", \"94\u00d727=\"),\n    @(\"60\u00d712=\", \"27\u00d791=\"),\n    @(\"83\u00d723=\", \"26\u00d752=\"),\n    @(\"16\u00d723=\", \"15\u00d723=\"),\n    @(\"11\u00d796=\", \"70\u00d798=\"),\n    @(\"38\u00d750=\", \"45\u00d739=\"),\n    @(\"83\u00d721=\", \"58\u00d798=\"),\n    @(\"47\u00d772=\", \"48\u00d785=\"),\n    @(\"74\u00d721=\", \"65\u00d732=\"),\n    @(\"99\u00d776=\", \"84\u00d739=\"),\n    @(\"58\u00d742=\", \"93\u00d720=\"),\n    @(\"93\u00d775=\", \"86\u00d744=\"),\n    @(\"30\u00d740=\", \"31\u00d755=\"),\n    @(\"78\u00d793=\", \"12\u00d778=\"),\n    @(\"90\u00d756=\", \"71\u00d724=\"),\n    @(\"81\u00d724=\", \"53\u00d758=\"),\n    @(\"58\u00d713=\", \"50\u00d761=\"),\n    @(\"88\u00d770=\", \"85\u00d760=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
